$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New values added to existing rows ---
$ws.Range("O3").Value = 40
$ws.Range("N6").Value = 15
$ws.Range("J7").Value = 15
$ws.Range("O7").Value = 30

# --- New rows 8, 9, 10 (inserted before the old summary row 13) ---
# Shared-string table order must match the authored workbook: the label
# text is registered in the order "Generator Kommentieren", "Java und
# JDBC", "Java und JDBC Kommentieren" - so write B9's text first even
# though its row is written after B8 numerically.
$ws.Range("B9").Value = "Generator Kommentieren"
$ws.Range("B8").Value = "Java und JDBC"
$ws.Range("B10").Value = "Java und JDBC Kommentieren"

$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 5

$ws.Range("O9").Value = 40

$ws.Range("O10").Value = 25

# --- Clear the old row 13 (sum row); its formulas move down to row 14 ---
$ws.Range("C13:R13").ClearContents()

# --- Row 14 becomes the new sum row, summing rows 3-11 (skipping row 8) ---
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($col in $cols) {
    $ws.Range($col + "14").Formula = "=" + $col + "3+" + $col + "4+" + $col + "5+" + $col + "6+" + $col + "7+" + $col + "9+" + $col + "10+" + $col + "11"
}
$ws.Range("R14").Formula = "=C14+D14+E14+F14+G14+H14+I14+J14+K14+L14+M14+N14+O14+P14+Q14"

# --- Row 15 (was row 14) divides row 14 by 60 ---
$cols2 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols2) {
    $ws.Range($col + "15").Formula = "=" + $col + "14/60"
}

# --- Update selection to match the authored state ---
$ws.Range("U8").Select()
